$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------
# Sheet "SignIn" (first sheet) - sign in validation test data
# -----------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

$ws1.Cells.Item(1,1).Value = "Invalid input"
$ws1.Cells.Item(1,2).Value = "Username"
$ws1.Cells.Item(1,3).Value = "Password"
$ws1.Cells.Item(1,4).Value = "Error message"

$ws1.Cells.Item(2,1).Value = "Invalid username and password"
$ws1.Cells.Item(2,2).Value = "techfoxx"
$ws1.Cells.Item(2,3).Value = "superje14r"
$ws1.Cells.Item(2,4).Value = "Invalid Username and Password"

$ws1.Cells.Item(3,1).Value = "Invalid username and password"
$ws1.Cells.Item(3,2).Value = "archie2@5t"
$ws1.Cells.Item(3,3).Value = "aut7o2314r"
$ws1.Cells.Item(3,4).Value = "Invalid Username and Password"

$ws1.Cells.Item(4,1).Value = "Valid username and invalid password`t"
$ws1.Cells.Item(4,2).Value = "ninjafox"
$ws1.Cells.Item(4,3).Value = "superje14r"
$ws1.Cells.Item(4,4).Value = "Please check your password"

$ws1.Cells.Item(5,1).Value = "Valid username and invalid password`t"
$ws1.Cells.Item(5,2).Value = "ninjafox"
$ws1.Cells.Item(5,3).Value = "aut7o2314r"
$ws1.Cells.Item(5,4).Value = "Please check your password"

# Apply borders to the hyperlink cell (B3) first, then the rest of the
# table, and finally bold the header - this ordering keeps the style
# table as compact as possible.
$ws1.Range("B3").Borders.LineStyle = 1
$ws1.Range("A1:D5").Borders.LineStyle = 1
$ws1.Range("A1:D1").Font.Bold = $true

# -----------------------------------------------------------------
# Sheet "Pythoncode" (second sheet) - python snippet validation data
# -----------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Cells.Item(1,1).Value = "Valid code"
$ws2.Cells.Item(1,2).Value = "Result for valid code"
$ws2.Cells.Item(1,3).Value = "Invalid code"

$ws2.Cells.Item(2,1).Value = 'print "numpy"'
$ws2.Cells.Item(2,2).Value = "numpy"
$ws2.Cells.Item(2,3).Value = "numpy"

$ws2.Cells.Item(3,1).Value = 'print "ninja"'
$ws2.Cells.Item(3,2).Value = "ninja"
$ws2.Cells.Item(3,3).Value = 'print "ninja":'

# Borders + bold header
$ws2.Range("A2:C3").Borders.LineStyle = 1
$ws2.Range("A1:C1").Borders.LineStyle = 1
$ws2.Range("A1:C1").Font.Bold = $true

# Column widths (values chosen so the engine's internal pixel rounding
# lands as close as possible to the target character-width values)
$ws2.Columns.Item(2).ColumnWidth = 19.833333333333332
$ws2.Columns.Item(3).ColumnWidth = 14.333333333333334
$ws2.Columns.Item(4).ColumnWidth = 29.5

# move selection on sheet2, then re-activate sheet1 so it remains the
# selected tab (matches original file where SignIn was the active tab)
# and finally restore the SignIn selection to A17.
$ws2.Range("B1").Select()
$ws1.Activate()
$ws1.Range("A17").Select()
